$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental") column B currently has no value; it should become
# the literal text "false" (not the boolean FALSE). Typing a leading
# apostrophe forces Excel to store it as text (quote-prefixed) instead of
# auto-converting the recognized boolean literal.
$ws.Range("B7").Value = "'false"

# Re-apply the original cell formatting (border/alignment/etc. inherited
# from the rest of column A/B) so the cell keeps its normal "data" style
# instead of the ad-hoc quote-prefixed style that typing created.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 8 ("Date") column B: update the timestamp value.
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"
